$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell G3 currently holds the placeholder text "PrintToPDF" instead of
# the real BSE notices/circulars URL. Fix the text and turn it into a proper
# hyperlink (matching the style of the other PDF_URL cells in the sheet).
$url = "https://www.bseindia.com/markets/MarketInfo/DispNewNoticesCirculars.aspx?page=20251121-72"

$cell = $ws.Range("G3")
$cell.Value = $url

$ws.Hyperlinks.Add($cell, $url) | Out-Null
$cell.Style = "Hyperlink"
